{"js": "// Remove the two runs of text \"asas\" and \"jasjas\" (concatenated as\n// \"asasjasjas\") that sit just before the \"_GoBack\" bookmark, leaving the\n// paragraph otherwise intact (pPr + bookmarkStart/bookmarkEnd remain).\nconst body = context.document.body;\n\nconst results = body.search(\"asasjasjas\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Delete the matched range's text/runs entirely (not just clearing the text),\n  // so the <w:r> elements themselves are removed from the XML.\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the stray runs of text \"asas\" and \"jasjas\" (they sit adjacent to\n# one another, forming \"asasjasjas\") that precede the \"_GoBack\" bookmark.\n# The paragraph's formatting (pPr) and the bookmark itself are left intact.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"asasjasjas\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
